$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.286.01'
$ws.Cells.Item(2, 5).Value = '  +0.18%  '

$ws.Cells.Item(3, 4).Value = '3.504.27'
$ws.Cells.Item(3, 5).Value = '  -0.56%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '''589.44'
$ws.Cells.Item(5, 5).Value = '  +0.34%  '

$ws.Cells.Item(6, 4).Value = '''134.29'
$ws.Cells.Item(6, 5).Value = '  +0.13%  '

$ws.Cells.Item(7, 5).Value = '  -0.02%  '

$ws.Cells.Item(8, 4).Value = '''0.488'
$ws.Cells.Item(8, 5).Value = '  -0.45%  '

$ws.Cells.Item(9, 5).Value = '  +0.12%  '

$ws.Cells.Item(10, 4).Value = '''7.31'
$ws.Cells.Item(10, 5).Value = '  +2.30%  '

$ws.Cells.Item(11, 5).Value = '  +2.30%  '

$ws.Cells.Item(12, 4).Value = '4.100.90'
$ws.Cells.Item(12, 5).Value = '  -0.59%  '

$ws.Cells.Item(13, 5).Value = '  +1.22%  '

$ws.Cells.Item(14, 5).Value = '  +0.97%  '

$ws.Cells.Item(15, 4).Value = '3.505.47'
$ws.Cells.Item(15, 5).Value = '  -0.58%  '

$ws.Cells.Item(16, 4).Value = '64.311.19'
$ws.Cells.Item(16, 5).Value = '  +0.16%  '

$ws.Cells.Item(17, 4).Value = '''25.69'
$ws.Cells.Item(17, 5).Value = '  -6.53%  '

$ws.Cells.Item(18, 4).Value = '''9.86'
$ws.Cells.Item(18, 5).Value = '  +0.54%  '

$ws.Cells.Item(19, 4).Value = '''5.76'
$ws.Cells.Item(19, 5).Value = '  +2.51%  '

$ws.Cells.Item(20, 4).Value = '''13.53'
$ws.Cells.Item(20, 5).Value = '  -2.78%  '

$ws.Cells.Item(21, 4).Value = '''393.14'
$ws.Cells.Item(21, 5).Value = '  +2.70%  '

$ws.Cells.Item(22, 5).Value = '  -0.11%  '

$ws.Cells.Item(23, 4).Value = '3.643.79'
$ws.Cells.Item(23, 5).Value = '  -0.62%  '

$ws.Cells.Item(24, 4).Value = '''74.65'
$ws.Cells.Item(24, 5).Value = '  +0.86%  '

$ws.Cells.Item(25, 4).Value = '''0.999'
$ws.Cells.Item(25, 5).Value = '  -0.03%  '

$ws.Cells.Item(26, 5).Value = '  -0.18%  '

$ws.Cells.Item(27, 5).Value = '  +0.07%  '

$ws.Cells.Item(28, 4).Value = '''7.35'
$ws.Cells.Item(28, 5).Value = '  -1.61%  '

$ws.Cells.Item(29, 5).Value = '  +0.80%  '

$ws.Cells.Item(30, 5).Value = '  -2.53%  '

$ws.Cells.Item(31, 5).Value = '  -7.75%  '

$ws.Cells.Item(32, 4).Value = '3.526.32'
$ws.Cells.Item(32, 5).Value = '  -0.34%  '

$ws.Cells.Item(33, 4).Value = '''0.153'
$ws.Cells.Item(33, 5).Value = '  +5.40%  '

$ws.Cells.Item(34, 5).Value = '  +0.04%  '

$ws.Cells.Item(35, 4).Value = '''23.45'
$ws.Cells.Item(35, 5).Value = '  -0.67%  '

$ws.Cells.Item(36, 4).Value = '''5.13'
$ws.Cells.Item(36, 5).Value = '  -5.20%  '

$ws.Cells.Item(37, 4).Value = '''6.88'
$ws.Cells.Item(37, 5).Value = '  -1.02%  '

$ws.Cells.Item(38, 2).Value = 'ImmutableX'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(38, 4).Value = '''1.55'
$ws.Cells.Item(38, 5).Value = '  -0.79%  '

$ws.Cells.Item(39, 2).Value = 'Monero'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(39, 4).Value = '''167.41'
$ws.Cells.Item(39, 5).Value = '  +4.31%  '

$ws.Cells.Item(40, 4).Value = '''0.0781'
$ws.Cells.Item(40, 5).Value = '  -0.80%  '

$ws.Cells.Item(41, 4).Value = '''0.811'
$ws.Cells.Item(41, 5).Value = '  -0.32%  '

$ws.Cells.Item(42, 5).Value = '  -0.03%  '

$ws.Cells.Item(43, 5).Value = '  -5.35%  '

$ws.Cells.Item(44, 5).Value = '  -0.47%  '

$ws.Cells.Item(45, 4).Value = '''1.66'
$ws.Cells.Item(45, 5).Value = '  +2.84%  '

$ws.Cells.Item(46, 4).Value = '''1.17'
$ws.Cells.Item(46, 5).Value = '  -4.16%  '

$ws.Cells.Item(47, 4).Value = '''6.76'
$ws.Cells.Item(47, 5).Value = '  -0.74%  '

$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).Value = '2.337.43'
$ws.Cells.Item(48, 5).Value = '  -5.61%  '

$ws.Cells.Item(49, 2).Value = 'SuiNetwork'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(49, 4).Value = '''0.892'
$ws.Cells.Item(49, 5).Value = '  -2.06%  '

$ws.Cells.Item(50, 5).Value = '  -1.49%  '

$ws.Cells.Item(51, 4).Value = '''21.14'
$ws.Cells.Item(51, 5).Value = '  -1.44%  '

